# Regenerate the K column (column G) using the new Strikeout-derived values
# instead of the old "Strike#" based values, per the commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 4
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 2
    13 = 2
    14 = 2
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 1
    21 = 2
    23 = 1
    24 = 1
    25 = 2
    26 = 0
    27 = 2
    28 = 3
    29 = 2
    30 = 1
    31 = 3
    32 = 1
    33 = 2
    34 = 1
    35 = 3
    36 = 3
    37 = 0
    38 = 1
    39 = 1
    40 = 0
    41 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
